# "Update code tinh luong % format cac bang"
#
#  - Sheet "Đơn sale chính": the wide per-order detail block (columns
#    G..AA: nhóm dịch vụ, sale chính, các cột thanh toán/bác sĩ/phụ phẫu,
#    tỉ lệ & tiền chiết khấu sale phụ, ...) is collapsed down to a short
#    G..N block (tên dịch vụ, đơn giá gốc, sale phụ, upsale, đơn giá, đã
#    thanh toán, tỉ lệ chiết khấu sale chính, chiết khấu sale chính).
#  - Sheet "Lương": relabel "Danh mục" -> "Danh mục lương" and refresh the
#    ngày công / phụ cấp / LONG XUYÊN salary figures that depend on it.

$wb = $excel.ActiveWorkbook

$wsOrders = $wb.Worksheets.Item("Đơn sale chính")
$wsSalary = $wb.Worksheets.Item("Lương")

# ---------------------------------------------------------------------
# 1) "Đơn sale chính" — header row.
# ---------------------------------------------------------------------
$wsOrders.Range("G1").Value = "Tên dịch vụ"
$wsOrders.Range("H1").Value = "Đơn giá gốc"
$wsOrders.Range("I1").Value = "Sale phụ"
$wsOrders.Range("J1").Value = "Upsale"
$wsOrders.Range("K1").Value = "Đơn giá"
$wsOrders.Range("L1").Value = "Đã thanh toán"
$wsOrders.Range("M1").Value = "Tỉ lệ chiết khấu sale chính"
$wsOrders.Range("N1").Value = "Chiết khấu sale chính"

# ---------------------------------------------------------------------
# "Đơn sale chính" — data row 2 (HD-LUXURY / Tiêm môi order).
# ---------------------------------------------------------------------
$wsOrders.Range("G2").Value = "Tiêm môi"
$wsOrders.Range("H2").Value = 1800000
$wsOrders.Range("I2").ClearContents()
$wsOrders.Range("J2").ClearContents()
$wsOrders.Range("K2").Value = 1800000
$wsOrders.Range("L2").Value = 1800000
$wsOrders.Range("M2").Value = 0.1
$wsOrders.Range("N2").Value = 180000

# ---------------------------------------------------------------------
# "Đơn sale chính" — total row 3. (G3 / I3 were already blank and stay
# blank, so they are left untouched.)
# ---------------------------------------------------------------------
$wsOrders.Range("H3").Value = 1800000
$wsOrders.Range("J3").Value = 0
$wsOrders.Range("K3").Value = 1800000
$wsOrders.Range("L3").Value = 1800000
$wsOrders.Range("M3").Value = 0.1
$wsOrders.Range("N3").Value = 180000

# Drop the now-unused tail columns (old N..AA) so the sheet's used range
# shrinks back down to A1:N3.
$wsOrders.Range("O1:AA3").ClearContents()

# ---------------------------------------------------------------------
# 2) "Lương" sheet — label + recomputed totals.
# ---------------------------------------------------------------------
$wsSalary.Range("A1").Value = "Danh mục lương"
$wsSalary.Range("B2").Value = 18
$wsSalary.Range("B3").Value = 630000
$wsSalary.Range("B12").Value = 2571428.571428571
$wsSalary.Range("B29").Value = 3381428.571428571
$wsSalary.Range("B31").Value = 3381428.571428571
